# Applies the "Generated Update from Main Repository" edit:
#   - Slide 1  ("Stuff" / Content Placeholder 2): full rewrite of the body
#     text to the new weekly-update bullets, including updated indent
#     levels and resetting autofit (normAutofit with no lnSpcReduction).
#   - Slide 18 ("Do we Select?" / Content Placeholder 2): the first
#     paragraph's three runs get collapsed back into a single run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1: Content Placeholder 2 - replace body text entirely.
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$shape1 = $s1.Shapes.Item("Content Placeholder 2")
$tf1 = $shape1.TextFrame
$tr1 = $tf1.TextRange

# Write a throwaway placeholder first so the engine doesn't try to splice the
# new text onto the old runs (it reuses formatting for overlapping substrings
# otherwise), then assign the real, final text in one shot (paragraphs are
# separated with carriage returns, same as real PowerPoint TextRange.Text).
$tr1.Text = "x"
$tr1.Text = "Feature selection basics (chapter 4)`rI added a review of multicollinearity to it, if you have old version, pull. `rMostly mechanics of different ways to select features. `rThis part is pretty straightforward, again make sure the general ‘make a pipe’ stuff is OK. `rAnd then… Next time, support vector machines (maybe some this time, if you’re fast). `rDimensionality and visualizing data. `rMargins and hinge loss. `rThen Natural Language Processing (NLP) – our first really cool topic. `rThe dimensionality bit is important for understanding here. "

# Paragraph indent levels (IndentLevel is 1-based: 1 -> <a:pPr lvl="0"/> i.e.
# no pPr element at all, 2 -> lvl="1", 3 -> lvl="2").
$tr1.Paragraphs(2, 1).IndentLevel = 2
$tr1.Paragraphs(3, 1).IndentLevel = 2
$tr1.Paragraphs(4, 1).IndentLevel = 2
$tr1.Paragraphs(6, 1).IndentLevel = 2
$tr1.Paragraphs(7, 1).IndentLevel = 2
$tr1.Paragraphs(8, 1).IndentLevel = 2
$tr1.Paragraphs(9, 1).IndentLevel = 3

# Paragraph 4 ("This part is pretty straightforward...") ends with a
# separately-formatted trailing fragment in the target deck, so split it into
# two runs by re-writing the tail of the paragraph through a Characters
# sub-range (this produces a second <a:r> without touching bold/italic/etc).
$para4 = $tr1.Paragraphs(4, 1)
$tail4 = $para4.Characters(80, 13)
$tail4.Text = "stuff is OK. "

# Restore the "shrink text on overflow" autofit with no active reduction yet
# (matches <a:normAutofit/> with no lnSpcReduction attribute in the target).
$tf1.AutoSize = 2

# ---------------------------------------------------------------------------
# Slide 18: Content Placeholder 2 - merge first paragraph's 3 runs into 1.
# ---------------------------------------------------------------------------
$s18 = $p.Slides.Item(18)
$shape18 = $s18.Shapes.Item("Content Placeholder 2")
$tr18 = $shape18.TextFrame.TextRange
$para1_18 = $tr18.Paragraphs(1, 1)

# Same placeholder trick: avoids partial substring run-splicing so the whole
# paragraph collapses back down to a single run instead of keeping the
# original three-run split.
$para1_18.Text = "x"
$para1_18.Text = "As noted, removing features makes the model more ‘efficient’ but doesn’t normally raise the ceiling. "
